# Update the "Metadata" worksheet (sheet 1) to reflect the new IG publication
# metadata: version bump, status change, new date, updated contact info, a new
# Jurisdiction row, and the resulting shift of the remaining metadata rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 12 for the "Jurisdiction" property. This pushes
# Description/Purpose/Copyright/Immutable down by one row (old rows 12-15
# become new rows 13-16). Use Insert() then copy the formatting from the row
# above (row 11, which already carries the correct border/alignment style) so
# the new row matches the existing formatting instead of Excel's default.
$ws.Rows.Item(12).Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Version: 0.1.6 -> 0.1.7
$ws.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: refresh publication timestamp
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# Contact details: first contact row now holds the organization contact
# (instead of the generic "No display for ContactDetail" placeholder), and
# the second (newly distinct) contact row holds the named contact person.
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# New row 12: Jurisdiction property with an empty value.
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

Write-Host "Metadata sheet updated"
